$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set E7, E24, E31, E32, E34, E35, E36, E54, E55, E67, E68 to "II"
$rowsII = @(7,24,31,32,34,35,36,54,55,67,68)
foreach ($r in $rowsII) {
    $ws.Range("E$r").Value = "II"
}

# Row 33: E33 = "II", F33 = "X"
$ws.Range("E33").Value = "II"
$ws.Range("F33").Value = "X"

# Update selection to E7, and remove topLeftCell by scrolling to A1
$ws.Range("E7").Select()

$wb.Save()
